$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.064.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.650.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5195'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2639'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06327'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07680'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.594'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.651.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.878.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5581'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8127'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.070.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.623'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '191.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.913'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1185'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.205'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05460'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.268'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.440'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.343'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.552'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.787'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9452'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01577'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.846'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.025.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8268'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.02'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.789.72'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈111'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9993'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.993'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4334'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05166'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.70%  '
